$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the timestamp column (H2:H51) to reflect the new ingestion run time.
$newTimestamp = "2025-04-05 21:48:36"

for ($row = 2; $row -le 51; $row++) {
    $ws.Cells.Item($row, 8).Value = $newTimestamp
}
